$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New guest rows (5-16), name in column A, table ("mesa") number in column B
$data = @(
    @("lucas", 1),
    @("maria", 3),
    @("natalia", 3),
    @("manuaela", 4),
    @("manuela", 5),
    @("pedro", 2),
    @("joao", 1),
    @("ronal", 3),
    @("francisco", 4),
    @("junior ", 5),
    @("felipe", 5),
    @("patricia", 6)
)

$row = 5
foreach ($entry in $data) {
    if ($row -eq 10) {
        # A10 previously held an empty, styled (underlined) cell; clear the
        # formatting before writing the real value so no style carries over.
        $ws.Range("A10").ClearFormats()
    }
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $row++
}

# Update the active selection to B20, as recorded in the saved workbook
$ws.Range("B20").Select()
